$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$urlBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fd9579c8128bb7c36c118e306cd579c7181ce517/e2e/"
$file1 = "23f6bd92-b7d1-4908-94ed-2075f84ea54d.md"
$file2 = "2fe9b83d-e653-41fc-9bec-d0f2cbcbf83b.md"

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet mirrors the Status text in columns E (zh-cn) / F (de-de) ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet (report row 2 = file1, row 3 = file2) ---
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Range("I2").Value = $file1
$wsZh.Range("J2").Value = "23f6bd92-b7d1-4908-94ed-2075f84ea54d.e9ab5f197ccc10b9c3db6e169080ed05ded0cdae.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-27 11:01:14"

$wsZh.Range("I3").Value = $file2
$wsZh.Range("J3").Value = "2fe9b83d-e653-41fc-9bec-d0f2cbcbf83b.9324d7989e3346b39dcea40cb541305b6d28540b.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-27 11:01:14"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($urlBase + $file1), "", "", $file1)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($urlBase + $file2), "", "", $file2)

# --- de-de sheet (report row 2 = file1, row 3 = file2) ---
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Range("I2").Value = $file1
$wsDe.Range("J2").Value = "23f6bd92-b7d1-4908-94ed-2075f84ea54d.e9ab5f197ccc10b9c3db6e169080ed05ded0cdae.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-27 11:01:20"

$wsDe.Range("I3").Value = $file2
$wsDe.Range("J3").Value = "2fe9b83d-e653-41fc-9bec-d0f2cbcbf83b.9324d7989e3346b39dcea40cb541305b6d28540b.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-27 11:01:20"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($urlBase + $file1), "", "", $file1)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($urlBase + $file2), "", "", $file2)

# --- Column width adjustments to reflect the new, longer content ---
# (ColumnWidth assignments land on a fixed pixel grid, so pick inputs whose
# rounded result matches the report's generated widths as closely as possible)
$wsOverview.Range("E1").ColumnWidth = 29.15
$wsOverview.Range("F1").ColumnWidth = 29.15

$wsZh.Range("C1").ColumnWidth = 29.15
$wsZh.Range("I1").ColumnWidth = 39.15
$wsZh.Range("J1").ColumnWidth = 39.15

$wsDe.Range("C1").ColumnWidth = 29.15
$wsDe.Range("I1").ColumnWidth = 39.15
$wsDe.Range("J1").ColumnWidth = 39.15
